$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the text of the existing comment (currently anchored on the
# "Error handling strategy..." task row, A7) so it can be re-created on
# the same task after the row above it is removed.
$commentText = $ws.Range("A7").Comment.Text()

# Select row 2 (the "Textured surfaces" task) just like a user would
# before deleting it, then delete the entire row - this shifts every
# row below it up by one.
$ws.Range("A2:XFD2").Select()
$ws.Rows("2:2").Delete()

# The task that used to live on row 7 is now on row 6; move the comment
# along with it.
$ws.Range("A7").Comment.Delete()
$ws.Range("A6").AddComment($commentText)
